$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1975.75
$ws.Range("J29").Value = 3877
$ws.Range("L29").Value = 11631
$ws.Range("N29").Value = -12193

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H58").Value = 10005
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 10005
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 30015
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -30315

$ws.Range("H87").Value = 93569.336

$ws.Range("H90").Value = 93569.336

$ws.Range("H98").Value = 582.2
$ws.Range("I98").Value = 582.2
$ws.Range("K98").Value = 582.2
$ws.Range("M98").Value = 915.8

$ws.Range("H112").Value = 1364.8077
$ws.Range("J112").Value = 1749.1666
$ws.Range("L112").Value = 5247.4998
$ws.Range("N112").Value = -7463.4998

$ws.Range("H122").Value = 582.2
$ws.Range("I122").Value = 582.2
$ws.Range("K122").Value = 1746.6
$ws.Range("M122").Value = 703.3999999999999

$ws.Range("H132").Value = 1254.3334
$ws.Range("I132").Value = 1254.3334
$ws.Range("K132").Value = 3763.0002
$ws.Range("M132").Value = -1233.0002

$ws.Range("H137").Value = 4133
$ws.Range("I137").Value = 1618.625
$ws.Range("K137").Value = 4855.875
$ws.Range("M137").Value = -2305.875

$ws.Range("H138").Value = 7504.073
$ws.Range("J138").Value = 11391.16
$ws.Range("L138").Value = 34173.48
$ws.Range("N138").Value = -44453.48

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 825.4286
$ws.Range("J2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("N2").Value = -1226

$ws.Range("H5").Value = 236
$ws.Range("I5").Value = 93.333336
$ws.Range("J5").Value = 450
$ws.Range("K5").Value = 93.333336
$ws.Range("L5").Value = 450
$ws.Range("M5").Value = 18.666664
$ws.Range("N5").Value = -674

$ws.Range("H32").Value = 3679.147
$ws.Range("I32").Value = 3293.6128
$ws.Range("K32").Value = 3293.6128
$ws.Range("M32").Value = -3006.6128

$ws.Range("H63").Value = 6000
$ws.Range("J63").Value = 6000
$ws.Range("L63").Value = 6000
$ws.Range("N63").Value = -7372

$ws.Range("H66").Value = 6000
$ws.Range("J66").Value = 6000
$ws.Range("L66").Value = 30000
$ws.Range("N66").Value = -36864

$ws.Range("H102").Value = 2189.2
$ws.Range("I102").Value = 1543.6666
$ws.Range("K102").Value = 1543.6666
$ws.Range("M102").Value = 78.33339999999998

$ws.Range("H116").Value = 825.4286
$ws.Range("J116").Value = 1000
$ws.Range("L116").Value = 1000
$ws.Range("N116").Value = -5588

$ws.Range("H122").Value = 1389.875
$ws.Range("I122").Value = 1389.875
$ws.Range("K122").Value = 4169.625
$ws.Range("M122").Value = -1719.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 825.4286
$ws.Range("J3").Value = 1000
$ws.Range("L3").Value = 1000
$ws.Range("N3").Value = -1228

$ws.Range("H4").Value = 236
$ws.Range("I4").Value = 93.333336
$ws.Range("J4").Value = 450
$ws.Range("K4").Value = 93.333336
$ws.Range("L4").Value = 450
$ws.Range("M4").Value = 21.666664
$ws.Range("N4").Value = -680

$ws.Range("H82").Value = 23977.857
$ws.Range("I82").Value = 8371
$ws.Range("J82").Value = 62995
$ws.Range("K82").Value = 8371
$ws.Range("L82").Value = 62995
$ws.Range("M82").Value = -7988
$ws.Range("N82").Value = -63761

$ws.Range("H85").Value = 23977.857
$ws.Range("I85").Value = 8371
$ws.Range("J85").Value = 62995
$ws.Range("K85").Value = 8371
$ws.Range("L85").Value = 62995
$ws.Range("M85").Value = -7045
$ws.Range("N85").Value = -65647

$ws.Range("H134").Value = 2139.8572
$ws.Range("I134").Value = 1946.2
$ws.Range("K134").Value = 5838.6
$ws.Range("M134").Value = -3303.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4983.2856
$ws.Range("I62").Value = 4972
$ws.Range("J62").Value = 4998.3335
$ws.Range("K62").Value = 4972
$ws.Range("L62").Value = 4998.3335
$ws.Range("M62").Value = -4348
$ws.Range("N62").Value = -6246.3335

$ws.Range("H65").Value = 4983.2856
$ws.Range("I65").Value = 4972
$ws.Range("J65").Value = 4998.3335
$ws.Range("K65").Value = 24860
$ws.Range("L65").Value = 24991.6675
$ws.Range("M65").Value = -21740
$ws.Range("N65").Value = -31231.6675

$ws.Range("H131").Value = 94900
$ws.Range("J131").Value = 94900
$ws.Range("L131").Value = 94900
$ws.Range("N131").Value = -104980

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1999
$ws.Range("I3").Value = 1999
$ws.Range("K3").Value = 5997
$ws.Range("M3").Value = -5885

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 89
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 45
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 45
$ws.Range("M2").Value = 13
$ws.Range("N2").Value = -271

$ws.Range("H18").Value = 20952
$ws.Range("J18").Value = 20952
$ws.Range("L18").Value = 20952
$ws.Range("N18").Value = -21538

$ws.Range("H43").Value = 20838
$ws.Range("J43").Value = 30757
$ws.Range("L43").Value = 30757
$ws.Range("N43").Value = -31059

$ws.Range("H46").Value = 34923.332
$ws.Range("J46").Value = 34923.332
$ws.Range("L46").Value = 34923.332
$ws.Range("N46").Value = -35235.332

$ws.Range("H57").Value = 14990
$ws.Range("J57").Value = 14990
$ws.Range("L57").Value = 14990
$ws.Range("N57").Value = -16630

$ws.Range("H80").Value = 5232
$ws.Range("I80").Value = 5672.6665
$ws.Range("J80").Value = 4901.5
$ws.Range("K80").Value = 5672.6665
$ws.Range("L80").Value = 4901.5
$ws.Range("M80").Value = -4674.6665
$ws.Range("N80").Value = -6897.5

$ws.Range("H83").Value = 5232
$ws.Range("I83").Value = 5672.6665
$ws.Range("J83").Value = 4901.5
$ws.Range("K83").Value = 28363.3325
$ws.Range("L83").Value = 24507.5
$ws.Range("M83").Value = -23371.3325
$ws.Range("N83").Value = -34491.5

$ws.Range("H102").Value = 3558.7
$ws.Range("I102").Value = 3085.875
$ws.Range("K102").Value = 3085.875
$ws.Range("M102").Value = -1463.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 40077
$ws.Range("J54").Value = 40077
$ws.Range("L54").Value = 40077
$ws.Range("N54").Value = -41117

$ws.Range("H81").Value = 6213.7
$ws.Range("I81").Value = 2304.7144
$ws.Range("K81").Value = 4609.4288
$ws.Range("M81").Value = -3548.4288

$ws.Range("H84").Value = 6213.7
$ws.Range("I84").Value = 2304.7144
$ws.Range("K84").Value = 23047.144
$ws.Range("M84").Value = -17743.144
